# OMOGENEIA_FOREIGN_COUNTRY_SECONDARY.docx
#
# Commit: "Modify protocol field in school transports' templates"
#
# The protocol-number line in the header table currently reads:
#     Αρ. Πρωτ.: ${protocol}
# It must read:
#     Αρ. Πρωτ.: Φ.15.1/${protocol}
#
# i.e. the literal text "Φ.15.1/" is inserted right before the merge
# field opener "${" that precedes the "protocol" placeholder.

$d = $word.ActiveDocument

# Single-quoted literal: PowerShell does not expand $-variables inside
# single quotes, so "${" is safe to embed here.
$needle = 'Αρ. Πρωτ.: ${'

# Locate the unique "Αρ. Πρωτ.: ${" run of text in the document body.
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if ($found) {
    # Collapse to the end of the match, then pull the start back over the
    # trailing "${" so we can insert immediately in front of it (i.e.
    # right after "Αρ. Πρωτ.: " and before "${protocol}").
    $rng.Collapse(0)
    $null = $rng.MoveStart(1, -2)
    $rng.InsertBefore("Φ.15.1/")
}
